$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Reshape the source data from a single wide query (A1:F3 - one row per
#    person, one column per category) into a tall layout that supports
#    multiple queries (A1:C5 - one row per category, one column per person).
# ---------------------------------------------------------------------------

# Drop the old D:F columns (they held "Cat 2".."Cat 4" / extra values).
$ws.Range("D1:F3").Clear()

# Header row: category-column header + the two series names.
$ws.Range("A1").Value = "Sales"
$ws.Range("B1").Value = "Person 1"
$ws.Range("C1").Value = "Person 2"

# Cat 0 row
$ws.Range("A2").Value = "Cat 0"
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 13

# Cat 1 row (formulas, same as before)
$ws.Range("A3").Value = "Cat 1"
$ws.Range("B3").Formula = "= MOD(B2+4,20)"
$ws.Range("C3").Formula = "= MOD(C2+4,20)"

# Cat 2 row
$ws.Range("A4").Value = "Cat 2"
$ws.Range("B4").Value = 18
$ws.Range("C4").Formula = "= MOD(B4+4,20)"

# Cat 3 row
$ws.Range("A5").Value = "Cat 3"
$ws.Range("B5").Value = 23
$ws.Range("C5").Formula = "= MOD(B5+4,20)"

# ---------------------------------------------------------------------------
# 2. Point the chart's two series (bar = Person 1, line = Person 2) at the
#    new layout: name cell in row 1, categories in A2:A5, values in the
#    person's own column.
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(BarLineOverlay!`$B`$1,BarLineOverlay!`$A`$2:`$A`$5,BarLineOverlay!`$B`$2:`$B`$5,1)"

$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(BarLineOverlay!`$C`$1,BarLineOverlay!`$A`$2:`$A`$5,BarLineOverlay!`$C`$2:`$C`$5,2)"

# No auto title on the chart any more.
$chart.HasTitle = $false

# Category axis picks up an explicit "General" source-linked number format.
$catAx = $chart.Axes(1)
$catAx.TickLabels.NumberFormatLinked = 1

# ---------------------------------------------------------------------------
# 3. Resize / reposition the chart on the sheet.
# ---------------------------------------------------------------------------
$co.Left = 239.0
$co.Top = 91.5
$co.Width = 433.0625
$co.Height = 201.0
